$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update the Date value (row 8, column B)
$ws.Range("B8").Value = "2023-02-01T09:05:11-06:00"

# Add the "true" value for the Experimental row (row 7, column B) as literal text,
# not an auto-converted Boolean. Write it via a formula producing the text "true",
# then convert the formula result to a static value (Paste Special - Values),
# which preserves Excel's literal-text cell type (t="s") instead of coercing
# the string "true" to a Boolean the way a direct .Value assignment would.
$ws.Range("B7").Formula = "=""true"""
$ws.Range("B7").Copy()
$ws.Range("B7").PasteSpecial(-4163)
